# Refresh the cryptos list "Price" (D) and "Volume(1h)" (E) columns with
# the latest scraped figures, per commit:
#   "Updated cryptos list on Sat Apr 20 13:00:10 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column strings are written with a leading apostrophe so Excel
# stores them as literal text: some look numeric ("1.00", "0.517", ...) and
# some use "." as a thousands separator ("63.869.52"), so a plain numeric
# write would reformat/round them. (In a single-quoted PowerShell string a
# literal ' is written as '' - e.g. '''1.00' is the 5-character string
# '1.00 - that leading apostrophe.) Style is reset back to Normal afterwards
# so no stray text-format style is left applied to the cell.

$ws.Range("D2").Value = '''63.869.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.76%  '
$ws.Range("D3").Value = '''3.062.52'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.24%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").Value = '''557.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.62%  '
$ws.Range("D6").Value = '''142.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.86%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '''3.062.07'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.00%  '
$ws.Range("D9").Value = '''0.517'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.29%  '
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("D11").Value = '''6.16'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.98%  '
$ws.Range("D12").Value = '''0.481'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.82%  '
$ws.Range("D13").Value = '''0.0000231'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.60%  '
$ws.Range("D14").Value = '''35.24'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").Value = '''3.564.15'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.33%  '
$ws.Range("D16").Value = '''63.898.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.90%  '
$ws.Range("D17").Value = '''3.066.07'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.38%  '
$ws.Range("E18").Value = '  +0.41%  '
$ws.Range("D19").Value = '''6.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").Value = '''489.17'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.56%  '
$ws.Range("D21").Value = '''14.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.62%  '
$ws.Range("D22").Value = '''0.687'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.05%  '
$ws.Range("D23").Value = '''14.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.59%  '
$ws.Range("E24").Value = '  -0.22%  '
$ws.Range("D25").Value = '''82.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("E27").Value = '  +0.44%  '
$ws.Range("D28").Value = '''8.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.84%  '
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("D30").Value = '''1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("D31").Value = '''26.61'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.74%  '
$ws.Range("E32").Value = '  -0.24%  '
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("D34").Value = '''5.69'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.75%  '
$ws.Range("D35").Value = '''6.22'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.80%  '
$ws.Range("D36").Value = '''55.22'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.60%  '
$ws.Range("D37").Value = '''0.0411'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").Value = '''443.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.27%  '
$ws.Range("D39").Value = '''0.0815'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.28%  '
$ws.Range("D40").Value = '''3.035.43'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.14%  '
$ws.Range("E41").Value = '  -5.77%  '
$ws.Range("D42").Value = '''8.34'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.07%  '
$ws.Range("E43").Value = '  +0.33%  '
$ws.Range("D44").Value = '''0.273'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.70%  '
$ws.Range("D45").Value = '''27.83'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.78%  '
$ws.Range("E46").Value = '  +2.70%  '
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("D49").Value = '''118.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.65%  '
$ws.Range("D50").Value = '''0.0₃0516'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.70%  '
$ws.Range("E51").Value = '  +2.14%  '
